$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 100
$ws.Range("C3").Value = 55
$ws.Range("C5").Value = 60

$ws.Range("C6").Select()
